# Weekly re-shuffle of the daily "Pepino dulce" price rows.
# Each destination row (2-43) receives the Fecha/Calidad/Volumen/Precio
# block that used to live on a different source row. L and M always mirror
# K (min = max = weighted avg for this market), and P is the derived
# $/Kg figure (round(K / Q)), so they travel together with K.
#
# Row 21 keeps its own original data (it is a fixed point of the permutation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (source row's original values land on
# the destination row)
$rowMap = @{
    2  = 11
    3  = 10
    4  = 3
    5  = 28
    6  = 9
    7  = 26
    8  = 38
    9  = 29
    10 = 20
    11 = 36
    12 = 40
    13 = 31
    14 = 6
    15 = 12
    16 = 24
    17 = 5
    18 = 8
    19 = 25
    20 = 42
    22 = 23
    23 = 32
    24 = 34
    25 = 19
    26 = 15
    27 = 14
    28 = 13
    29 = 35
    30 = 17
    31 = 22
    32 = 37
    33 = 16
    34 = 41
    35 = 33
    36 = 39
    37 = 43
    38 = 4
    39 = 7
    40 = 30
    41 = 27
    42 = 18
    43 = 2
}

# 1) Snapshot the columns that move (D, I, J, K, L, M, P) for every data
#    row *before* writing anything, so overlapping permutation cycles
#    don't clobber a value before it has been read.
$snapshot = @{}
for ($r = 2; $r -le 43; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value()
        I = $ws.Cells.Item($r, 9).Value()
        J = $ws.Cells.Item($r, 10).Value()
        K = $ws.Cells.Item($r, 11).Value()
        L = $ws.Cells.Item($r, 12).Value()
        M = $ws.Cells.Item($r, 13).Value()
        P = $ws.Cells.Item($r, 16).Value()
    }
}

# 2) Write each destination row's block from the snapshot of its mapped
#    source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $src = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value  = $src.D
    $ws.Cells.Item($destRow, 9).Value  = $src.I
    $ws.Cells.Item($destRow, 10).Value = $src.J
    $ws.Cells.Item($destRow, 11).Value = $src.K
    $ws.Cells.Item($destRow, 12).Value = $src.L
    $ws.Cells.Item($destRow, 13).Value = $src.M
    $ws.Cells.Item($destRow, 16).Value = $src.P
}
